# Swap the order of "Recorded By" entries in column G so that
# values written as "<email>, System" become "System, <email>".
# Only touches cells that actually match this exact two-part pattern;
# all other cells (single value, different ordering, 3+ parts, etc.)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string]) {
        $parts = $val -split ', '
        if ($parts.Count -eq 2 -and $parts[1] -eq 'System' -and $parts[0] -ne 'System') {
            $cell.Value2 = 'System, ' + $parts[0]
        }
    }
}
